$d = $word.ActiveDocument

# 1. Delete the "Meta description" paragraph that currently follows the
#    title/Heading1 paragraph ("Play Fruits Slot Free - Colorful Slot
#    Machine Game"). Locate it by its leading label text rather than a
#    bare index so the script stays correct even if earlier content
#    shifts.
$metaIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Meta description")) {
    $metaIdx = $i
    break
  }
}
if ($metaIdx -gt 0) {
  $d.Paragraphs.Item($metaIdx).Range.Delete()
}

# 2. Insert a new bold paragraph with the title text right before the
#    final paragraph (the one that holds the italic image-prompt text),
#    mirroring the exact run layout used elsewhere in the doc
#    (an empty leading run followed by the formatted text run).
#    First carve out a brand-new, empty paragraph immediately before the
#    final one (this leaves all surrounding paragraphs untouched), then
#    fill that empty paragraph's range with the desired run XML.
$promptIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Create a feature image")) {
    $promptIdx = $i
    break
  }
}
if ($promptIdx -le 0) {
  $promptIdx = $d.Paragraphs.Count
}

$lastPara = $d.Paragraphs.Item($promptIdx)
$insPoint = $lastPara.Range.Duplicate
$insPoint.Collapse(1)
$insPoint.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($promptIdx)
$newRange = $newPara.Range
$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fruits Slot Free - Colorful Slot Machine Game</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$newRange.InsertXML($xmlFrag)

# 3. Replace the final paragraph's text (the old image-generation prompt)
#    with the meta-description copy, keeping its italic formatting intact.
$oldText = "Create a feature image that showcases the fun and quirky vibe of " + [char]34 + "Fruits" + [char]34 + ". The image should be in a cartoon style and feature a happy Maya warrior wearing glasses. The warrior should be surrounded by colorful fruity symbols while holding a mobile device that displays the game " + [char]34 + "Fruits" + [char]34 + ". The background should have a tropical setting with palm trees, blue skies, and vibrant fruit patterns. The overall design should be eye-catching and playful, representing the excitement and appeal of this mobile-first slot game."
$newText = "Read Fruits slot review. Play Fruits online slot game for free. Fun, colorful slot machine with a classic feel and modern twist."
[void]$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
